$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 edits ---
# B3: was a numeric phone (502992932); becomes empty.
$ws.Range("B3").ClearContents()

# J3: last_activity date text changes from 2025-11-20 to 2025-11-18.
# Force text storage (avoid auto date-serial conversion), then drop the
# extra number-format style so the cell keeps the default style.
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2025-11-18"
$ws.Range("J3").Style = "Normal"

# --- Row 6 edit ---
# J6: last_activity date text changes from 2025-11-27 to 2025-11-17.
$ws.Range("J6").NumberFormat = "@"
$ws.Range("J6").Value = "2025-11-17"
$ws.Range("J6").Style = "Normal"

# --- Row 7 edit ---
# B7: was a numeric phone (503535395); becomes the text "0503535395".
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "0503535395"
$ws.Range("B7").Style = "Normal"

# --- Remove trailing rows 8-11 (duplicate/garbage entries) ---
$ws.Rows("8:11").Delete()
